# Replace project with updated version:
# Append 5 new daily rows (221-225) to the portfolio values table on Sheet1,
# extending the data range from A1:F220 to A1:F225.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: date(serial), cash, buying_power, portfolio_value, SPX USD, FX_EURUSD
$newRows = @(
    @{ Row = 221; A = 46062; B = 2110.53; C = 115759.3613751221; D = 115759.3613751221; E = 6964.81982421875;  F = 1.182578206062317 },
    @{ Row = 222; A = 46063; B = 2110.53; C = 115368.84;         D = 115368.84;         E = 6941.81005859375;  F = 1.190405368804932 },
    @{ Row = 223; A = 46064; B = 2110.53; C = 116938.43;         D = 116938.43;         E = 6941.47021484375;  F = 1.188961744308472 },
    @{ Row = 224; A = 46065; B = 2110.53; C = 114662.9;          D = 114662.9;          E = 6832.759765625;    F = 1.18760621547699  },
    @{ Row = 225; A = 46066; B = 2110.53; C = 115932.09;         D = 115932.09;         E = 6836.169921875;    F = 1.186844944953918 }
)

$firstNewRow = $newRows[0].Row
$lastNewRow = $newRows[$newRows.Count - 1].Row

# Row 220 is the last existing data row; copy its cell formatting (date style on
# column A, plain/general style on B:F) down onto the new rows before filling values,
# so the new cells pick up the same style index (s="2" on column A) as the rest of
# the table rather than creating brand-new style entries.
$ws.Range("A220:F220").Copy() | Out-Null
$ws.Range("A${firstNewRow}:F${lastNewRow}").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $r.A
    $ws.Range("B$row").Value2 = $r.B
    $ws.Range("C$row").Value2 = $r.C
    $ws.Range("D$row").Value2 = $r.D
    $ws.Range("E$row").Value2 = $r.E
    $ws.Range("F$row").Value2 = $r.F
}
